$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 with the new regression transaction data (Debt Case
# Management local execution merged with regression Feb02), replacing
# the prior "dayfivemo" run's values.
$ws.Range("A7").Value = "01-25-2019 17:17:10"
$ws.Range("B7").Value = "FT19012500004"

# C7 holds a 17-digit numeric string; entering it as a plain value would
# let Excel auto-coerce it to a double and lose precision, so enter it
# the way a user would to force text ('-prefix), then drop the resulting
# quote-prefix formatting so the cell is left with no explicit style,
# matching how the rest of the sheet's text cells are stored.
$ws.Range("C7").Value = "'20190125044036506"
$ws.Range("C7").ClearFormats()

$ws.Range("G7").Value = "Regression daysixtw"
